# Update (Removed Auto Arima)
# Fill in previously-blank/auto-arima-affected forecast columns
# (Amazon Mean / P70 / P80 / P90 Forecast) on the "Forecast Comparison"
# sheet with the recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$values = @{
    2  = @(39, 48, 56, 70)
    3  = @(36, 45, 56, 74)
    4  = @(28, 35, 42, 52)
    5  = @(28, 34, 41, 52)
    6  = @(27, 33, 41, 55)
    7  = @(27, 33, 42, 58)
    8  = @(27, 33, 43, 59)
    9  = @(28, 34, 45, 64)
    10 = @(27, 33, 42, 57)
    11 = @(27, 33, 44, 63)
    12 = @(27, 32, 43, 62)
    13 = @(29, 35, 46, 66)
    14 = @(27, 33, 44, 63)
    15 = @(27, 32, 44, 64)
    16 = @(27, 33, 44, 64)
    17 = @(26, 31, 42, 61)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Cells.Item($row, 4).Value = $rowValues[0]  # D: Amazon Mean Forecast
    $ws.Cells.Item($row, 5).Value = $rowValues[1]  # E: Amazon P70 Forecast
    $ws.Cells.Item($row, 6).Value = $rowValues[2]  # F: Amazon P80 Forecast
    $ws.Cells.Item($row, 7).Value = $rowValues[3]  # G: Amazon P90 Forecast
}
